$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 4

$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 8

$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 5
$ws.Range("I4").Value = 3

$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 3

$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 7
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 1

$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 10
$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 1

$ws.Range("K9").Select()

# The source row's column-F width was widened (content-driven best-fit
# recalculation in the original edit). Reproduce the resulting width as
# closely as this engine's ColumnWidth quantization allows.
$ws.Columns.Item(6).ColumnWidth = 18.8333333333

